# Generate Report for Handoff
#
# A new handoff just completed for the 8d0cffe5-57db-4a96-890b-cf2d0d0dd30d
# file, so its "Latest Handoff Datetime" cell (column E, row 6) is refreshed
# with the newer timestamp on both the zh-cn and de-de status sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E6").Value = "2016-03-18 22:33:02"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E6").Value = "2016-03-18 22:33:06"
